$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reuse the header formatting from the adjacent "sum" column (G1) for the
# new "Save" header in H1, so it picks up the existing bold/border/centered
# style rather than creating a brand-new style entry.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
